$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress percentages (C2, C3) from 90% to 100%
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1

# C8 was a text status ("en proceso"); it's now a completed numeric value (100)
$ws.Range("C8").Value = 100

# Mark rows 21 and 22 as "en proceso" (same status text used elsewhere, e.g. C4/C8/C17)
$ws.Range("C21").Value = $ws.Range("C17").Value2
$ws.Range("C22").Value = $ws.Range("C17").Value2

# Move the active selection to C9
$ws.Range("C9").Select()
